# Update coin Price / Volume(1h) figures per the Jan 12 2023 GitHub Actions symbol-list refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'285.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'2.81%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'28.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'4.50%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.041"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'3.72%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06483"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.02%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.221"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'3.85%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.332"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'12.93%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9128"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'4.29%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.14%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.06495"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'25.58%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07616"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.49%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.02985"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.99%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'-0.31%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'2.11%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006529"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'2.45%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006117"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.07%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.459"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.58%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.369"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.87%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-1.41%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'1.34%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1341"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-0.53%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'1.82%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.1556"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'3.73%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04460"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.97%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'0.49%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004321"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'11.65%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D28").Value = "'0.0001181"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'-9.13%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.0001637"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'-15.65%"
$ws.Range("E29").Style = "Normal"
$ws.Range("D40").Value = "'0.04148"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.39%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006701"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.96%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1232"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'5.05%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002112"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'7.24%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01190"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'3.68%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005394"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.66%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-0.06%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.041"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'20.97%"
$ws.Range("E47").Style = "Normal"
